# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = '@'
    $cell.Value = $text
    $cell.Style = 'Normal'
}

Set-TextValue 2 4 '44.435.92'
Set-TextValue 2 5 '  +0.69%  '
Set-TextValue 3 4 '2.246.90'
Set-TextValue 3 5 '  -0.29%  '
Set-TextValue 4 5 '  +0.25%  '
Set-TextValue 5 4 '306.39'
Set-TextValue 5 5 '  -0.24%  '
Set-TextValue 6 4 '93.48'
Set-TextValue 6 5 '  -5.43%  '
Set-TextValue 7 5 '  -0.58%  '
Set-TextValue 8 4 '1.00'
Set-TextValue 8 5 '  +0.24%  '
Set-TextValue 9 4 '0.524'
Set-TextValue 9 5 '  -2.14%  '
Set-TextValue 10 4 '34.65'
Set-TextValue 10 5 '  -2.88%  '
Set-TextValue 11 4 '0.0812'
Set-TextValue 11 5 '  -1.43%  '
Set-TextValue 12 4 '7.16'
Set-TextValue 12 5 '  -2.30%  '
Set-TextValue 13 5 '  -0.16%  '
Set-TextValue 14 2 'Polygon'
Set-TextValue 14 3 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 14 4 '0.837'
Set-TextValue 14 5 '  -0.42%  '
Set-TextValue 15 2 'WrappedEther'
Set-TextValue 15 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 15 4 '2.233.51'
Set-TextValue 15 5 '  -0.92%  '
Set-TextValue 16 4 '13.56'
Set-TextValue 16 5 '  -2.34%  '
Set-TextValue 17 4 '44.133.24'
Set-TextValue 17 5 '  +0.31%  '
Set-TextValue 18 4 '0.0₃0962'
Set-TextValue 18 5 '  -1.56%  '
Set-TextValue 19 4 '12.35'
Set-TextValue 19 5 '  -3.51%  '
Set-TextValue 20 4 '6.36'
Set-TextValue 20 5 '  +0.10%  '
Set-TextValue 21 4 '65.67'
Set-TextValue 21 5 '  +0.36%  '
Set-TextValue 22 2 'BitcoinCash'
Set-TextValue 22 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 22 4 '238.25'
Set-TextValue 22 5 '  -1.31%  '
Set-TextValue 23 2 'PancakeSwap'
Set-TextValue 23 3 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 23 4 '2.95'
Set-TextValue 23 5 '  -0.10%  '
Set-TextValue 24 4 '1.99'
Set-TextValue 24 5 '  +0.33%  '
Set-TextValue 25 5 '  +0.16%  '
Set-TextValue 26 4 '38.51'
Set-TextValue 26 5 '  +3.35%  '
Set-TextValue 27 4 '2.21'
Set-TextValue 27 5 '  +3.22%  '
Set-TextValue 28 4 '9.82'
Set-TextValue 28 5 '  -3.13%  '
Set-TextValue 29 4 '5.95'
Set-TextValue 29 5 '  -3.75%  '
Set-TextValue 30 4 '20.09'
Set-TextValue 30 5 '  -0.12%  '
Set-TextValue 31 4 '152.96'
Set-TextValue 31 5 '  -2.73%  '
Set-TextValue 32 4 '0.0801'
Set-TextValue 32 5 '  -2.97%  '
Set-TextValue 33 4 '2.66'
Set-TextValue 33 5 '  +0.13%  '
Set-TextValue 34 5 '  -12.21%  '
Set-TextValue 35 5 '  +1.46%  '
Set-TextValue 36 5 '  +0.51%  '
Set-TextValue 37 4 '1.80'
Set-TextValue 37 5 '  -3.58%  '
Set-TextValue 38 4 '3.47'
Set-TextValue 38 5 '  +2.43%  '
Set-TextValue 39 4 '14.75'
Set-TextValue 39 5 '  -3.90%  '
Set-TextValue 40 4 '3.83'
Set-TextValue 40 5 '  -1.66%  '
Set-TextValue 41 4 '0.0301'
Set-TextValue 41 5 '  -1.69%  '
Set-TextValue 43 4 '1.735.05'
Set-TextValue 43 5 '  -1.31%  '
Set-TextValue 44 4 '80.69'
Set-TextValue 44 5 '  -7.61%  '
Set-TextValue 45 4 '0.192'
Set-TextValue 45 5 '  -0.43%  '
Set-TextValue 46 4 '99.78'
Set-TextValue 46 5 '  -1.87%  '
Set-TextValue 47 4 '1.62'
Set-TextValue 47 5 '  +4.51%  '
Set-TextValue 48 4 '4.93'
Set-TextValue 48 5 '  -4.21%  '
Set-TextValue 49 4 '8.18'
Set-TextValue 49 5 '  -0.69%  '
Set-TextValue 50 4 '55.41'
Set-TextValue 50 5 '  -0.17%  '
Set-TextValue 51 4 '69.03'
Set-TextValue 51 5 '  -1.93%  '
